$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add 'hole_id' header to A1, copying the formatting used by the other
# header cells (bold font, border, centered alignment) from B1.
$ws1.Range("B1").Copy()
$ws1.Range("A1").PasteSpecial(-4122)
$ws1.Range("A1").Value = "hole_id"

# Replace the numeric index values in column A (rows 2-29) with the
# corresponding hole_id string values.
$holeIds = @(
  "BRG_16_02",
  "BRG_05_09",
  "BRG_01_02",
  "BRG_05_11",
  "ECO_09_01",
  "ECO_09_02",
  "BRG_05_14",
  "BRG_16_03",
  "BRG_05_02",
  "BRG_01_05",
  "BRG_01_07",
  "BRG_01_04",
  "BRG_16_07",
  "BRG_01_06",
  "BRG_13_02",
  "BRG_05_15",
  "ECO_09_05",
  "BRG_08_01",
  "BRG_16_09",
  "BRG_05_12",
  "BRG_01_03",
  "BRG_01_08",
  "ECO_09_04",
  "BRG_05_03",
  "BRG_05_13",
  "BRG_13_01",
  "BRG_01_09",
  "BRG_05_10"
)

$row = 2
foreach ($holeId in $holeIds) {
  $ws1.Cells.Item($row, 1).Value = $holeId
  $row = $row + 1
}
